$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in Group column for the "Reusable Solutions" row (A3)
$ws.Range("A3").Value = "Reusable Solutions"

# Update row heights (rows 2-5)
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 120
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 45

# Update column widths (A and B) to new slightly narrower widths
$ws.Columns.Item(1).ColumnWidth = 14.022135416666666
$ws.Columns.Item(2).ColumnWidth = 152.59244791666666
